$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.256.44'
$ws.Range('E2').Value = '  +1.07%  '

$ws.Range('D3').Value = '3.565.72'
$ws.Range('E3').Value = '  +4.73%  '

$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').Value = '607.12'
$ws.Range('E5').Value = '  +1.94%  '

$ws.Range('D6').Value = '145.31'
$ws.Range('E6').Value = '  +2.20%  '

$ws.Range('D7').Value = '3.565.79'
$ws.Range('E7').Value = '  +4.90%  '

$ws.Range('E8').Value = '  +0.18%  '

$ws.Range('E9').Value = '  +3.47%  '

$ws.Range('E10').Value = '  +1.48%  '

$ws.Range('E11').Value = '  +1.06%  '

$ws.Range('E12').Value = '  +1.55%  '

$ws.Range('D13').Value = '4.172.17'
$ws.Range('E13').Value = '  +4.74%  '

$ws.Range('E14').Value = '  +3.86%  '

$ws.Range('E15').Value = '  +1.21%  '

$ws.Range('D16').Value = '3.562.12'
$ws.Range('E16').Value = '  +4.88%  '

$ws.Range('D17').Value = '66.372.47'
$ws.Range('E17').Value = '  +1.36%  '

$ws.Range('E18').Value = '  -0.78%  '

$ws.Range('D19').Value = '11.54'
$ws.Range('E19').Value = '  +11.25%  '

$ws.Range('E20').Value = '  +1.68%  '

$ws.Range('E21').Value = '  +1.15%  '

$ws.Range('D22').Value = '432.12'
$ws.Range('E22').Value = '  +3.66%  '

$ws.Range('E23').Value = '  +4.93%  '

$ws.Range('D24').Value = '78.80'
$ws.Range('E24').Value = '  +1.86%  '

$ws.Range('D25').Value = '3.707.94'
$ws.Range('E25').Value = '  +4.47%  '

$ws.Range('E26').Value = '  -0.01%  '

$ws.Range('E27').Value = '  +7.33%  '

$ws.Range('D28').Value = '8.05'
$ws.Range('E28').Value = '  +2.69%  '

$ws.Range('E29').Value = '  +3.94%  '

$ws.Range('D30').Value = '9.21'
$ws.Range('E30').Value = '  -0.64%  '

$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  -0.01%  '

$ws.Range('E32').Value = '  +0.38%  '

$ws.Range('D33').Value = '0.159'
$ws.Range('E33').Value = '  -0.65%  '

$ws.Range('D34').Value = '3.560.36'
$ws.Range('E34').Value = '  +4.70%  '

$ws.Range('E35').Value = '  +3.84%  '

$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '1.76'
$ws.Range('E36').Value = '  +3.38%  '

$ws.Range('B37').Value = 'USDe'
$ws.Range('C37').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.05%  '

$ws.Range('E38').Value = '  +4.35%  '

$ws.Range('E39').Value = '  +1.63%  '

$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  +0.05%  '

$ws.Range('D41').Value = '171.55'
$ws.Range('E41').Value = '  +1.30%  '

$ws.Range('E42').Value = '  -0.97%  '

$ws.Range('D43').Value = '5.23'
$ws.Range('E43').Value = '  +3.41%  '

$ws.Range('E44').Value = '  +3.25%  '

$ws.Range('E45').Value = '  +1.99%  '

$ws.Range('D46').Value = '46.08'
$ws.Range('E46').Value = '  +1.29%  '

$ws.Range('E47').Value = '  +3.73%  '

$ws.Range('D48').Value = '26.06'
$ws.Range('E48').Value = '  -2.84%  '

$ws.Range('E49').Value = '  +3.57%  '

$ws.Range('D50').Value = '7.14'
$ws.Range('E50').Value = '  +0.53%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '23.29'
$ws.Range('E51').Value = '  +14.11%  '
